# Apply the row permutation described by the diff.
# Each physical row (2-13) receives the full record previously
# held by a different row; only the cells whose value actually
# changes are written (plus clearing cells that must become blank).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 12
$ws.Range("A2").Value = 111815518
$ws.Range("B2").Value = 77515
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 458250.901553072
$ws.Range("R2").Value = 7054618.376188213

# Row 3 <- original row 2
$ws.Range("A3").Value = 111815516
$ws.Range("B3").Value = 89423
$ws.Range("E3").Value = 5432
$ws.Range("F3").Value = "Granticka"
$ws.Range("G3").Value = "Porodaedalea chrysoloma"
$ws.Range("H3").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q3").Value = 458289.5512131723
$ws.Range("R3").Value = 7054475.069158822
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("AC3").ClearContents()

# Row 4 <- original row 11
$ws.Range("A4").Value = 111815508
$ws.Range("B4").Value = 56398
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("Q4").Value = 458162.4570845839
$ws.Range("R4").Value = 7054329.489790585
$ws.Range("AC4").Value = "ringhack"
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

# Row 5 <- original row 9
$ws.Range("A5").Value = 111815513
$ws.Range("Q5").Value = 458173.7327805056
$ws.Range("R5").Value = 7054711.474791372
$ws.Range("AC5").Value = "ringhack gamla"

# Row 6 <- original row 3
$ws.Range("A6").Value = 111815507
$ws.Range("B6").Value = 56398
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("Q6").Value = 458151.5539710881
$ws.Range("R6").Value = 7054482.225765129
$ws.Range("AC6").Value = "ringhack gamla"
$ws.Range("K6").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()

# Row 7 <- original row 4
$ws.Range("A7").Value = 111815515
$ws.Range("B7").Value = 89423
$ws.Range("E7").Value = 5432
$ws.Range("F7").Value = "Granticka"
$ws.Range("G7").Value = "Porodaedalea chrysoloma"
$ws.Range("H7").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q7").Value = 458161.9437607233
$ws.Range("R7").Value = 7054459.400503729
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("AC7").ClearContents()

# Row 8 <- original row 6
$ws.Range("A8").Value = 111815519
$ws.Range("B8").Value = 77515
$ws.Range("E8").Value = 6425
$ws.Range("F8").Value = "Garnlav"
$ws.Range("G8").Value = "Alectoria sarmentosa"
$ws.Range("H8").Value = "(Ach.) Ach."
$ws.Range("Q8").Value = 458215.7474518137
$ws.Range("R8").Value = 7054621.063481365

# Row 9 <- original row 13
$ws.Range("A9").Value = 111815509
$ws.Range("Q9").Value = 458176.2590895323
$ws.Range("R9").Value = 7054362.673967168

# Row 10 <- original row 7
$ws.Range("A10").Value = 111815512
$ws.Range("B10").Value = 56398
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("Q10").Value = 458154.6107204149
$ws.Range("R10").Value = 7054646.336103803
$ws.Range("AC10").Value = "ringhack"
$ws.Range("K10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()

# Row 11 <- original row 8
$ws.Range("A11").Value = 111815514
$ws.Range("B11").Value = 89423
$ws.Range("E11").Value = 5432
$ws.Range("F11").Value = "Granticka"
$ws.Range("G11").Value = "Porodaedalea chrysoloma"
$ws.Range("H11").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q11").Value = 458153.7808649908
$ws.Range("R11").Value = 7054482.19637617
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("AC11").ClearContents()

# Row 12 <- original row 5
$ws.Range("A12").Value = 111815510
$ws.Range("B12").Value = 56398
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = "Tretåig hackspett"
$ws.Range("G12").Value = "Picoides tridactylus"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("Q12").Value = 458203.7272220219
$ws.Range("R12").Value = 7054385.000644128
$ws.Range("AC12").Value = "ringhack"
$ws.Range("K12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

# Row 13 <- original row 10
$ws.Range("A13").Value = 111815517
$ws.Range("B13").Value = 77515
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("Q13").Value = 458250.8216980004
$ws.Range("R13").Value = 7054375.482693202
$ws.Range("K13").ClearContents()
$ws.Range("L13").ClearContents()
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("AC13").ClearContents()
